# word-style-reference-ttmanu.docx update:
#   - Heading2/3/4: tighten paragraph spacing (before/after) and switch to
#     double line spacing (line=480, lineRule=auto) for an "auth names"
#     10pt-friendly rhythm.
#   - Heading6 ("auth names" / institutional affiliations style): set an
#     explicit 10pt run size.
#   - Bibliography: give it an explicit Arial typeface.

$d = $word.ActiveDocument

# --- Heading 2 --------------------------------------------------------
$h2 = $d.Styles("Heading2")
$h2.ParagraphFormat.SpaceBefore     = 14   # 280 twips
$h2.ParagraphFormat.SpaceAfter      = 0
$h2.ParagraphFormat.LineSpacingRule = 5    # wdLineSpaceMultiple
$h2.ParagraphFormat.LineSpacing     = 24   # -> w:line="480" w:lineRule="auto"

# --- Heading 3 --------------------------------------------------------
$h3 = $d.Styles("Heading3")
$h3.ParagraphFormat.SpaceBefore     = 12   # 240 twips
$h3.ParagraphFormat.SpaceAfter      = 0
$h3.ParagraphFormat.LineSpacingRule = 5
$h3.ParagraphFormat.LineSpacing     = 24

# --- Heading 4 --------------------------------------------------------
$h4 = $d.Styles("Heading4")
$h4.ParagraphFormat.SpaceBefore     = 12   # 240 twips
$h4.ParagraphFormat.SpaceAfter      = 0
$h4.ParagraphFormat.LineSpacingRule = 5
$h4.ParagraphFormat.LineSpacing     = 24

# --- Heading 6 (author / institutional-affiliation line) --------------
$h6 = $d.Styles("Heading6")
$h6.Font.Size = 10

# --- Bibliography -------------------------------------------------------
$bib = $d.Styles("Bibliography")
$bib.Font.Name = "Arial"
